$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" date column (C) for rows 2-7 from 45182 (2023-09-13)
# to 45184 (2023-09-15), keeping the existing date number format/style.
$ws.Range("C2:C7").Value = 45184
